$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("F2").Value = 25.8200000000006
$ws.Range("H2").Value = 0.02301193295792536
$ws.Range("I2").Value = 0.02301193295792536
$ws.Range("L2").Value = 7.060592151256505
$ws.Range("M2").Value = "[0.7273602421680092, 13.393824060345]"
$ws.Range("N2").Value = 0.0296996412860957
$ws.Range("O2").Value = 0.0296996412860957
$ws.Range("P2").Value = -1.371105502467618
$ws.Range("Q2").Value = "[-2.578684660604235, -0.16352634433100066]"
$ws.Range("R2").Value = 0.02696124826939617
$ws.Range("S2").Value = 0.02696124826939617
$ws.Range("T2").Value = 9.979234250362786
$ws.Range("U2").Value = "[6.323523951203326, 13.634944549522247]"
$ws.Range("V2").Value = 0.000001724877121356272
$ws.Range("W2").Value = 0.000001724877121356272
$ws.Range("X2").Value = 5.634394394394526
$ws.Range("Y2").Value = 0.6719919919920097
$ws.Range("Z2").Value = 10.59679679679704
$ws.Range("F3").Value = 25.8200000000006
$ws.Range("H3").Value = 0.2635806660112034
$ws.Range("I3").Value = 0.2635806660112034
$ws.Range("L3").Value = 4.487445569449578
$ws.Range("M3").Value = "[-2.6156169502062205, 11.590508089105377]"
$ws.Range("N3").Value = 0.2097545925011284
$ws.Range("O3").Value = 0.2097545925011284
$ws.Range("P3").Value = -1.849105585896695
$ws.Range("Q3").Value = "[-4.9749745525315845, 1.2767633807381937]"
$ws.Range("R3").Value = 0.2397265129616635
$ws.Range("S3").Value = 0.2397265129616635
$ws.Range("T3").Value = 9.551191717273273
$ws.Range("U3").Value = "[5.615703162425007, 13.48668027212154]"
$ws.Range("V3").Value = 0.00001333959439442367
$ws.Range("W3").Value = 0.00001333959439442367
$ws.Range("X3").Value = 7.598678678678858
$ws.Range("Y3").Value = -5.246706706706822
$ws.Range("Z3").Value = 20.44406406406454
$ws.Range("F4").Value = 25.8200000000006
$ws.Range("H4").Value = 0.07826808567034627
$ws.Range("I4").Value = 0.07826808567034627
$ws.Range("L4").Value = 5.917558968587938
$ws.Range("M4").Value = "[-0.8755729792492701, 12.710690916425147]"
$ws.Range("N4").Value = 0.08615177569021548
$ws.Range("O4").Value = 0.08615177569021548
$ws.Range("P4").Value = -2.239053022378311
$ws.Range("Q4").Value = "[-4.723395561253121, 0.24528951649649944]"
$ws.Range("R4").Value = 0.07615437789616331
$ws.Range("S4").Value = 0.07615437789616331
$ws.Range("T4").Value = 10.57651371898055
$ws.Range("U4").Value = "[6.98556821551161, 14.167459222449484]"
$ws.Range("V4").Value = 0.0000003936386139269388
$ws.Range("W4").Value = 0.0000003936386139269388
$ws.Range("X4").Value = 9.201121121121332
$ws.Range("Y4").Value = -1.00798798798801
$ws.Range("Z4").Value = 19.41023023023067
$ws.Range("F5").Value = 25.8200000000006
$ws.Range("H5").Value = 0.1143664286208084
$ws.Range("I5").Value = 0.1143664286208084
$ws.Range("L5").Value = 6.174490640934735
$ws.Range("M5").Value = "[-1.3464854709864813, 13.695466752855951]"
$ws.Range("N5").Value = 0.1051861796531293
$ws.Range("O5").Value = 0.1051861796531293
$ws.Range("P5").Value = -2.80510575275485
$ws.Range("Q5").Value = "[-4.679369237779392, -0.9308422677303079]"
$ws.Range("R5").Value = 0.004220257505968217
$ws.Range("S5").Value = 0.004220257505968217
$ws.Range("T5").Value = 10.45253156047906
$ws.Range("U5").Value = "[6.488048528448685, 14.417014592509426]"
$ws.Range("V5").Value = 0.00000325288913383126
$ws.Range("W5").Value = 0.00000325288913383126
$ws.Range("X5").Value = 11.52724724724751
$ws.Range("Y5").Value = 3.825185185185272
$ws.Range("Z5").Value = 19.22930930930976
$ws.Range("F6").Value = 25.8200000000006
$ws.Range("H6").Value = 0.1066468800218204
$ws.Range("I6").Value = 0.1066468800218204
$ws.Range("L6").Value = 6.478242326414515
$ws.Range("M6").Value = "[-0.6047889636103037, 13.561273616439333]"
$ws.Range("N6").Value = 0.07205174516993962
$ws.Range("O6").Value = 0.07205174516993962
$ws.Range("P6").Value = 3.08813211794312
$ws.Range("Q6").Value = "[1.578658170272349, 4.597606065613891]"
$ws.Range("R6").Value = 0.0001599876982683046
$ws.Range("S6").Value = 0.0001599876982683046
$ws.Range("T6").Value = 10.37260257610723
$ws.Range("U6").Value = "[6.2444876759470525, 14.500717476267408]"
$ws.Range("V6").Value = 0.000007511578462349178
$ws.Range("W6").Value = 0.000007511578462349178
$ws.Range("X6").Value = 13.12968968968999
$ws.Range("Y6").Value = 6.926686686686848
$ws.Range("Z6").Value = 19.33269269269314
$ws.Range("F7").Value = 25.8200000000006
$ws.Range("H7").Value = 0.1228062399297913
$ws.Range("I7").Value = 0.1228062399297913
$ws.Range("L7").Value = 5.468519059077748
$ws.Range("M7").Value = "[-1.8975995354215076, 12.834637653577003]"
$ws.Range("N7").Value = 0.1418305709880472
$ws.Range("O7").Value = 0.1418305709880472
$ws.Range("P7").Value = 2.547237286694427
$ws.Range("Q7").Value = "[-0.5346053564667317, 5.629079929855585]"
$ws.Range("R7").Value = 0.1029177240088897
$ws.Range("S7").Value = 0.1029177240088897
$ws.Range("T7").Value = 10.27759778912196
$ws.Range("U7").Value = "[6.526210812696688, 14.028984765547223]"
$ws.Range("V7").Value = 0.000001612042596965679
$ws.Range("W7").Value = 0.000001612042596965679
$ws.Range("X7").Value = 15.35243243243279
$ws.Range("Y7").Value = 2.68796796796803
$ws.Range("Z7").Value = 28.01689689689755
$ws.Range("F8").Value = 25.8200000000006
$ws.Range("H8").Value = 0.1868113990027263
$ws.Range("I8").Value = 0.1868113990027263
$ws.Range("L8").Value = 5.811492107120861
$ws.Range("M8").Value = "[-2.33385754459718, 13.956841758838902]"
$ws.Range("N8").Value = 0.1576322916770871
$ws.Range("O8").Value = 0.1576322916770871
$ws.Range("P8").Value = 2.786237328408966
$ws.Range("Q8").Value = "[-0.3018947895341544, 5.874369446352086]"
$ws.Range("R8").Value = 0.07584834250993411
$ws.Range("S8").Value = 0.07584834250993411
$ws.Range("T8").Value = 10.71103175304052
$ws.Range("U8").Value = "[6.315540378154505, 15.106523127926536]"
$ws.Range("V8").Value = 0.00001248798256159311
$ws.Range("W8").Value = 0.00001248798256159311
$ws.Range("X8").Value = 14.37029029029062
$ws.Range("Y8").Value = 1.679979979980018
$ws.Range("Z8").Value = 27.06060060060123
$ws.Range("F9").Value = 23.99000000000031
$ws.Range("H9").Value = 0.06668797728669928
$ws.Range("I9").Value = 0.06668797728669928
$ws.Range("L9").Value = 6.984884024885025
$ws.Range("M9").Value = "[-0.5974094374307253, 14.567177487200775]"
$ws.Range("N9").Value = 0.07009466974519607
$ws.Range("O9").Value = 0.07009466974519607
$ws.Range("P9").Value = 1.792500312859041
$ws.Range("Q9").Value = "[0.24528951649650033, 3.3397111092215823]"
$ws.Range("R9").Value = 0.02415399698355314
$ws.Range("S9").Value = 0.02415399698355314
$ws.Range("T9").Value = 11.35461564924714
$ws.Range("U9").Value = "[7.210502495115316, 15.498728803378954]"
$ws.Range("V9").Value = 0.000001609193148954802
$ws.Range("W9").Value = 0.000001609193148954802
$ws.Range("X9").Value = 17.14600600600623
$ws.Range("Y9").Value = 11.2385585585587
$ws.Range("Z9").Value = 23.05345345345376
$ws.Range("F10").Value = 23.99000000000031
$ws.Range("H10").Value = 0.1818293536365412
$ws.Range("I10").Value = 0.1818293536365412
$ws.Range("L10").Value = 5.191327734609423
$ws.Range("M10").Value = "[-2.2833590190068005, 12.666014488225647]"
$ws.Range("N10").Value = 0.16871657239579
$ws.Range("O10").Value = 0.16871657239579
$ws.Range("P10").Value = 2.207605648468504
$ws.Range("Q10").Value = "[-0.9308422677303092, 5.346053564667317]"
$ws.Range("R10").Value = 0.1634467022732435
$ws.Range("S10").Value = 0.1634467022732435
$ws.Range("T10").Value = 10.98032923845504
$ws.Range("U10").Value = "[7.056120159224125, 14.904538317685951]"
$ws.Range("V10").Value = 0.000001081286607140797
$ws.Range("W10").Value = 0.000001081286607140797
$ws.Range("X10").Value = 15.56108108108128
$ws.Range("Y10").Value = 3.57808808808813
$ws.Range("Z10").Value = 27.54407407407443
$ws.Range("B11").Value = 0
$ws.Range("F11").Value = 23.99000000000031
$ws.Range("H11").Value = 0.1244367734979475
$ws.Range("I11").Value = 0.1244367734979475
$ws.Range("L11").Value = 5.775580461348397
$ws.Range("M11").Value = "[-1.407654926940717, 12.958815849637512]"
$ws.Range("N11").Value = 0.1123468461362056
$ws.Range("O11").Value = 0.1123468461362056
$ws.Range("P11").Value = 2.270500396288119
$ws.Range("Q11").Value = "[-0.8553685703467702, 5.396369362923008]"
$ws.Range("R11").Value = 0.1504282306770115
$ws.Range("S11").Value = 0.1504282306770115
$ws.Range("T11").Value = 10.25143468552783
$ws.Range("U11").Value = "[6.315694095295411, 14.187175275760241]"
$ws.Range("V11").Value = 0.000004036528725048782
$ws.Range("W11").Value = 0.000004036528725048782
$ws.Range("X11").Value = 15.32094094094114
$ws.Range("Y11").Value = 3.385975975976018
$ws.Range("Z11").Value = 27.25590590590626
$ws.Range("F12").Value = 23.99000000000031
$ws.Range("H12").Value = 0.1301520961468657
$ws.Range("I12").Value = 0.1301520961468657
$ws.Range("L12").Value = 5.986945615245073
$ws.Range("M12").Value = "[-1.8854049133729731, 13.85929614386312]"
$ws.Range("N12").Value = 0.1325893333532211
$ws.Range("O12").Value = 0.1325893333532211
$ws.Range("P12").Value = 2.157289850212811
$ws.Range("Q12").Value = "[-0.6352369529781168, 4.949816653403738]"
$ws.Range("R12").Value = 0.1267280556429065
$ws.Range("S12").Value = 0.1267280556429065
$ws.Range("T12").Value = 10.64709123070423
$ws.Range("U12").Value = "[6.520972493306577, 14.773209968101884]"
$ws.Range("V12").Value = 0.000004757496128826588
$ws.Range("W12").Value = 0.000004757496128826588
$ws.Range("X12").Value = 15.7531931931934
$ws.Range("Y12").Value = 5.090970970971036
$ws.Range("Z12").Value = 26.41541541541576
$ws.Range("F13").Value = 23.99000000000031
$ws.Range("H13").Value = 0.2237227787488385
$ws.Range("I13").Value = 0.2237227787488385
$ws.Range("L13").Value = 5.236629459452352
$ws.Range("M13").Value = "[-3.221618362348101, 13.694877281252806]"
$ws.Range("N13").Value = 0.2188628361594034
$ws.Range("O13").Value = 0.2188628361594034
$ws.Range("P13").Value = 2.333395144107735
$ws.Range("Q13").Value = "[-0.8050527720910781, 5.471843060306547]"
$ws.Range("R13").Value = 0.1412557942207542
$ws.Range("S13").Value = 0.1412557942207542
$ws.Range("T13").Value = 10.60557121292213
$ws.Range("U13").Value = "[6.367730109271455, 14.843412316572795]"
$ws.Range("V13").Value = 0.000008038370278740103
$ws.Range("W13").Value = 0.000008038370278740103
$ws.Range("X13").Value = 15.080800800801
$ws.Range("Y13").Value = 3.097807807807847
$ws.Range("Z13").Value = 27.06379379379414
$ws.Range("F14").Value = 23.99000000000031
$ws.Range("H14").Value = 0.07345436366112312
$ws.Range("I14").Value = 0.07345436366112312
$ws.Range("L14").Value = 6.278139329476804
$ws.Range("M14").Value = "[-0.8946665304791921, 13.450945189432801]"
$ws.Range("N14").Value = 0.08471042772936399
$ws.Range("O14").Value = 0.08471042772936399
$ws.Range("P14").Value = 2.19502669890458
$ws.Range("Q14").Value = "[-0.4276842851733855, 4.817737682982545]"
$ws.Range("R14").Value = 0.09878356796979082
$ws.Range("S14").Value = 0.09878356796979082
$ws.Range("T14").Value = 10.00139468793505
$ws.Range("U14").Value = "[6.2085539269003895, 13.794235448969719]"
$ws.Range("V14").Value = 0.000003244719856132505
$ws.Range("W14").Value = 0.000003244719856132505
$ws.Range("X14").Value = 15.60910910910931
$ws.Range("Y14").Value = 5.595265265265338
$ws.Range("Z14").Value = 25.62295295295328
